# Refresh the cryptocurrency price / 1h-volume snapshot (GitHub Actions bot run).
# Most cells are plain text swaps; a handful of Price cells (column D) now
# hold digit-only strings (e.g. "252.96") that Excel would otherwise auto-
# convert to a Number on assignment, so those are briefly forced to Text
# format, written, then had their format cleared again to match the
# original (unstyled) text cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.142.77'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').Value = '1.901.43'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('E4').Value = '  -0.39%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '252.96'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.87%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.696'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.43%  '
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.12'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.38%  '
$ws.Range('E9').Value = '  +2.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.82'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0751'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.77%  '
$ws.Range('E12').Value = '  -1.31%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '13.07'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +5.71%  '
$ws.Range('D14').Value = '2.178.08'
$ws.Range('E14').Value = '  -0.09%  '
$ws.Range('E15').Value = '  +3.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.97'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.83%  '
$ws.Range('D17').Value = '1.905.75'
$ws.Range('E17').Value = '  +0.16%  '
$ws.Range('D18').Value = '35.142.78'
$ws.Range('E18').Value = '  -0.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.55'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.66%  '
$ws.Range('E20').Value = '  +1.62%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '242.73'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.94'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.05'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +5.34%  '
$ws.Range('E24').Value = '  -0.32%  '
$ws.Range('E25').Value = '  +5.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.29'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.74'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.98%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.55'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.30%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.50'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.69%  '
$ws.Range('E30').Value = '  -0.94%  '
$ws.Range('D31').Value = '4.128.85'
$ws.Range('E31').Value = '  -0.30%  '
$ws.Range('B32').Value = 'WEMIXToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.04'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +13.62%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0606'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +6.19%  '
$ws.Range('E34').Value = '  +4.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.57'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +7.65%  '
$ws.Range('E36').Value = '  +2.46%  '
$ws.Range('E37').Value = '  -0.37%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.855'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -8.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.00'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '103.84'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +15.52%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '17.32'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +8.17%  '
$ws.Range('E42').Value = '  +2.39%  '
$ws.Range('E43').Value = '  +0.61%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0652'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.43'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.53%  '
$ws.Range('D46').Value = '1.320.01'
$ws.Range('E46').Value = '  -1.94%  '
$ws.Range('E47').Value = '  +0.11%  '
$ws.Range('E48').Value = '  -1.71%  '
$ws.Range('E49').Value = '  +0.93%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '11.94'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -6.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0746'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +5.46%  '
